# "Restored working slate for triggerLogoFlip"
# Update the RFLNUM values (column B) on Sheet1 back to their working values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 4

# Move the sheet's active selection to G9 (was D10).
$ws.Range("G9").Select()

# Restore the workbook window width (cosmetic view state captured in the
# saved file alongside the sheet selection).
$excel.ActiveWindow.Width = 14400
